# "updated legacy GSC export data"
#
# The GSC (Google Search Console) video-indexing export is a rolling
# date window. The refreshed export drops the oldest date row
# (2025-10-29) from the "Chart" sheet; every following row — which
# already held the next day's figures, all the way through the newest
# date (2026-01-26) that the previous export had already appended at
# the bottom — simply shifts up by one row. Deleting the top data row
# (row 2, just below the header) reproduces exactly that shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows.Item(2).Delete()
